$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 110
$ws.Range("F4").Value = 406
$ws.Range("F5").Value = 969
$ws.Range("F6").Value = 5305
$ws.Range("F7").Value = 451
$ws.Range("F8").Value = 637
$ws.Range("F9").Value = 916
$ws.Range("F10").Value = 818
$ws.Range("F17").Value = 1772
$ws.Range("F18").Value = 1454
$ws.Range("F19").Value = 835
$ws.Range("F22").Value = 309
$ws.Range("F23").Value = 517
$ws.Range("F24").Value = 137
$ws.Range("F25").Value = 1046
$ws.Range("F28").Value = 2619
$ws.Range("F32").Value = 92
$ws.Range("F33").Value = 26
$ws.Range("F34").Value = 279
$ws.Range("F40").Value = 646
$ws.Range("F41").Value = 83
$ws.Range("F42").Value = 48
$ws.Range("F43").Value = 45
$ws.Range("F44").Value = 61

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 150
$ws.Range("F6").Value = 105

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 234

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 234
$ws.Range("F4").Value = 110
$ws.Range("F5").Value = 969
$ws.Range("F7").Value = 5305
$ws.Range("F8").Value = 451
$ws.Range("F9").Value = 637
$ws.Range("F11").Value = 150
$ws.Range("F12").Value = 916
$ws.Range("F13").Value = 818
$ws.Range("F15").Value = 105
$ws.Range("F23").Value = 1772
$ws.Range("F24").Value = 1454
$ws.Range("F25").Value = 835
$ws.Range("F27").Value = 309
$ws.Range("F29").Value = 517
$ws.Range("F30").Value = 137
$ws.Range("F31").Value = 1046
$ws.Range("F33").Value = 2619
$ws.Range("F36").Value = 92
$ws.Range("F37").Value = 26
$ws.Range("F38").Value = 279
$ws.Range("F43").Value = 83
$ws.Range("F44").Value = 48
$ws.Range("F45").Value = 45
$ws.Range("F46").Value = 61
